$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.032120108604431
$ws.Range("B1").Value = 1.378272771835327
$ws.Range("C1").Value = 2.228050947189331
$ws.Range("D1").Value = 4.465785026550293
$ws.Range("E1").Value = 1.993577718734741
